$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.260.28"
$ws.Range("E2").Value = "  +4.87%  "
$ws.Range("D3").Value = "1.909.51"
$ws.Range("E3").Value = "  +5.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5130"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2965"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06774"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("D11").Value = "1.915.22"
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07373"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6934"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.866"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").Value = "30.269.04"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008092"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "2.162.02"
$ws.Range("E21").Value = "  +5.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.834"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.738"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.135"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "136.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.002"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.235"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08793"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05099"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.153"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7158"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.833"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9731"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01689"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.065"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4293"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.93%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.651"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05755"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.505"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3808"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.08%  "
